# Resident Services_Requirements.xlsx
# Adds a new "Reg Proc" column (20th column, letter T) to the "Table2"
# ListObject on the "Details" worksheet, and populates the new column plus
# one updated comment cell (S8) with clarification notes gathered during
# review.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")
$lo = $ws.ListObjects.Item("Table2")

# Grow the table by one column (A2:S16 -> A2:T16); Excel names it
# "Column20" until we overwrite the header cell below.
$newCol = $lo.ListColumns.Add()

# Give the new column a sensible width to match its siblings.
$ws.Columns.Item(20).ColumnWidth = 32.08984375

# --- Header (matches the other row-2 header cells' look) ---------------
$ws.Range("S2").Copy()
$ws.Range("T2").PasteSpecial(-4122)
$ws.Range("T2").Value = "Reg Proc"

# --- New column body, row by row (format first, then text) -------------
$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial(-4122)
$ws.Range("T5").Value = "No Mapping of such kind from Reg Processor`nID Repo- Might not be there in ID Repo as well"

$ws.Range("S6").Copy()
$ws.Range("T6").PasteSpecial(-4122)
$ws.Range("T6").Value = "ID Repo- need to know "

$ws.Range("Q9").Copy()
$ws.Range("T9").PasteSpecial(-4122)
$ws.Range("T9").Value = "Under processing`nProcessed"

$ws.Range("Q8").Copy()
$ws.Range("T8").PasteSpecial(-4122)
$ws.Range("T8").Value = "Under processing`nProcessed`n"

# --- Updated "Comments" cell on row 8 (column S) ------------------------
$ws.Range("Q8").Copy()
$ws.Range("S8").PasteSpecial(-4122)
$ws.Range("S8").Value = "Reg proc`nArchival policy"

$ws.Range("R10").Copy()
$ws.Range("T10").PasteSpecial(-4122)
$ws.Range("T10").Value = "E-UIN Generation"

$ws.Range("R7").Copy()
$ws.Range("T7").PasteSpecial(-4122)
$ws.Range("T7").Value = "there shud be a label as Res_Service`nReg Client packet needs to be understood`nService from Reg proc needs to be developed"

$ws.Range("S5").Copy()
$ws.Range("T4").PasteSpecial(-4122)
$ws.Range("T4").Value = "When UIN IS needed to be generated`n1.the Acknowledgment from Print queue- what needs to be done`nTime period `n2. If there is a print failure- no need to handle from MOSIP`nUser Story ?"

# --- View state: keep selection on the newly-edited cell ----------------
$ws.Range("T4").Select()
